# Add a new "Swiss" test-data sheet (Switzerland market) based on the
# existing "Czech" sheet, matching the layout/styles of the other market
# sheets but trimmed down to a single attached-functionality entry.

$wb = $excel.ActiveWorkbook

# Duplicate the Czech sheet (keeps styles, merges, column widths, etc.)
# and place the copy at the very end of the workbook.
$czech = $wb.Worksheets.Item("Czech")
$czech.Copy($null, $czech)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Swiss"

# Drop the two extra "Attached Functionality" rows that don't apply to
# the Swiss market (Transmission Unit + Fire Brigade Panel / + Keysafe).
$ws.Range("A11:A12").EntireRow.Delete()

# Market-specific text. Fill "Fire Brigade Panel - LocalIO" before the
# user-story code so the new shared-string entries land in the same
# order as the source workbook.
$ws.Range("B2").Value = "Switzerland Market"
$ws.Range("A10").Value = "Fire Brigade Panel - LocalIO"
$ws.Range("B4").Value = "NGC-3476/T2650/T2660"

# Widen column B to fit the new market label.
$ws.Columns("B").ColumnWidth = 21.8

# Leave the cursor where the author left it on the new sheet.
$ws.Range("B9").Select()

# The author had selected the whole Czech sheet (click on the corner /
# Ctrl+A) before switching away to build the new tab.
$czech.Cells.Select()

$ws.Activate()
